$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 420-421, pushing the existing rows 420:437 down to 422:439
$ws.Rows("420:421").Insert()

# New row 420 data (Calidad: Primera)
$ws.Cells.Item(420, 1).Value = 3
$ws.Cells.Item(420, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(420, 3).Value = "Coquimbo"
$ws.Cells.Item(420, 4).Value = 44509
$ws.Cells.Item(420, 5).Value = 5
$ws.Cells.Item(420, 6).Value = 100112023
$ws.Cells.Item(420, 7).Value = "Brócoli"
$ws.Cells.Item(420, 8).Value = "Sin especificar"
$ws.Cells.Item(420, 9).Value = "Primera"
$ws.Cells.Item(420, 10).Value = 3000
$ws.Cells.Item(420, 11).Value = 500
$ws.Cells.Item(420, 12).Value = 550
$ws.Cells.Item(420, 13).Value = 523
$ws.Cells.Item(420, 14).Value = "$/unidad"
$ws.Cells.Item(420, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(420, 16).Value = 523
$ws.Cells.Item(420, 17).Value = 1
$ws.Cells.Item(420, 18).Value = "Hortaliza"

# New row 421 data (Calidad: Segunda)
$ws.Cells.Item(421, 1).Value = 3
$ws.Cells.Item(421, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(421, 3).Value = "Coquimbo"
$ws.Cells.Item(421, 4).Value = 44509
$ws.Cells.Item(421, 5).Value = 5
$ws.Cells.Item(421, 6).Value = 100112023
$ws.Cells.Item(421, 7).Value = "Brócoli"
$ws.Cells.Item(421, 8).Value = "Sin especificar"
$ws.Cells.Item(421, 9).Value = "Segunda"
$ws.Cells.Item(421, 10).Value = 1500
$ws.Cells.Item(421, 11).Value = 400
$ws.Cells.Item(421, 12).Value = 400
$ws.Cells.Item(421, 13).Value = 400
$ws.Cells.Item(421, 14).Value = "$/unidad"
$ws.Cells.Item(421, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(421, 16).Value = 400
$ws.Cells.Item(421, 17).Value = 1
$ws.Cells.Item(421, 18).Value = "Hortaliza"

# Ensure the date cells keep the date/time number format used by the rest of column D
$ws.Range("D420:D421").NumberFormat = $ws.Cells.Item(422, 4).NumberFormat
